# Applies the "Gaussian Quadrature Scheme" update to alpha4F-HW10.xlsx
#   1. Renames the worksheet tab from "alpha4F-HW10.xpc" to "alpha4F"
#   2. Refreshes a handful of row-13/row-15 values with re-computed
#      (last-bit) floating point results
#   3. Appends a new data row (row 16, HKL index 14 / "HexGrid-60degTilt5degRes")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet/tab
$ws.Name = "alpha4F"

# 2. Tiny recomputed values in the existing "11" (row 13) data row
$ws.Range("C13").Value = 0.9735639125594217
$ws.Range("F13").Value = 0.9735639125594217
$ws.Range("H13").Value = 1.001209575530364
$ws.Range("L13").Value = 0.9847077525664679
$ws.Range("M13").Value = 0.9904368462661726

# Tiny recomputed value in the "13" (row 15) data row
$ws.Range("H15").Value = 0.7463180020601147

# 3. New data row 16 -- HKL index 14, "HexGrid-60degTilt5degRes" (shared string 13)
$ws.Range("A16").Value = 14
$ws.Range("A13").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.623327862283947
$ws.Range("D16").Value = 0.4783232991450093
$ws.Range("E16").Value = 0.9523107647407856
$ws.Range("F16").Value = 1.623327862283947
$ws.Range("G16").Value = 0.6828143876837972
$ws.Range("H16").Value = 1.35554530157708
$ws.Range("I16").Value = 1.02730437444241
$ws.Range("J16").Value = 0.4783232991450093
$ws.Range("K16").Value = 0.7153170319428974
$ws.Range("L16").Value = 1.169322447113422
$ws.Range("M16").Value = 1.019937664978838
